$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data refresh cycles the Fecha/Volumen/Precio values through the
# existing rows (2-10). Only columns D, J, K, L, M, P change; everything
# else stays put.
$rows = @{
    2  = @{ D = 44907; J = 2300; K = 900;  L = 1000; M = 952;  P = 952 }
    3  = @{ D = 44883; J = 290;  K = 1400; L = 1500; M = 1434; P = 1434 }
    4  = @{ D = 44210; J = 1450; K = 1600; L = 1700; M = 1650; P = 1650 }
    5  = @{ D = 44175; J = 1400; K = 1900; L = 2000; M = 1950; P = 1950 }
    6  = @{ D = 44537; J = 800;  K = 1300; L = 1400; M = 1350; P = 1350 }
    7  = @{ D = 44895; J = 200;  K = 1200; L = 1300; M = 1255; P = 1255 }
    8  = @{ D = 44638; J = 800;  K = 2500; L = 2800; M = 2650; P = 2650 }
    9  = @{ D = 44893; J = 3300; K = 1200; L = 1300; M = 1261; P = 1261 }
    10 = @{ D = 44200; J = 1500; K = 1400; L = 1500; M = 1450; P = 1450 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("P$r").Value = $vals.P
}
